$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.50"
$ws.Range("D4").Value = "'5.309"
$ws.Range("D6").Value = "'3.379"
$ws.Range("D7").Value = "'6.379"
$ws.Range("D8").Value = "'0.8061"
$ws.Range("D9").Value = "'0.9500"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1426"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07443"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03182"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03085"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09242"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.572"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001620"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04724"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005823"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006371"
$ws.Range("D20").Value = "'0.004971"
$ws.Range("D21").Value = "'0.001045"
$ws.Range("D22").Value = "'0.0001504"
$ws.Range("D23").Value = "'0.0003106"
$ws.Range("E23").Value = "22UpBotsUBXT"
$ws.Range("D25").Value = "'2.098"
$ws.Range("D26").Value = "'0.3283"
$ws.Range("D40").Value = "'0.03942"
$ws.Range("D41").Value = "'0.006987"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").Value = "'0.1032"
$ws.Range("D43").Value = "'0.003161"
$ws.Range("D44").Value = "'0.008152"
$ws.Range("D45").Value = "'0.00005950"
$ws.Range("D46").Value = "'0.00000000752"
$ws.Range("D47").Value = "'0.0005513"
$ws.Range("E47").Value = "46ACDXExchangeACXTWorstin24h"
$ws.Range("D48").Value = "'0.6838"
$ws.Range("D49").Value = "'0.04035"
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("D51").Value = "'0.01012"

Write-Host "Applied all changes"
